$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Long text values (PT/EN course syllabus content) ---
$objetivosPt = "Desenvolver conceitos básicos da Estatística, com o apoio computacional, que permitam ao engenheiro trabalhar com o  fenômeno da aleatoriedade presente nos diversos campos de conhecimento da engenharia."
$docenteNome = "4894221 - Mariana Pereira de Melo"
$programaResumidoPt = "Estatística Descritiva, Modelos de Probabilidade, Teorema Central do Limite, Intervalos de Confiança, Testes de Hipóteses, ANOVA, Modelos de Regressão Linear."
$shortSyllabusEn = "Descriptive statistics, Probability models, Central limit theorem, Confidence intervals, Hypothesis test, ANOVA, Linear regression models."
$programaPt = "1)Estatística Descritiva: População e amostra; apresentação gráfica dos dados; medidas de posição; medidas de dispersão.2)Amostragem: Amostragem aleatória simples com reposição; amostragem aleatória simples sem reposição.3)Conceitos de Probabilidade: Conceitos básicos de probabilidade; operações com eventos; probabilidade condicional; independência; Teorema de Bayes.4)Variáveis Aleatórias discretas: Caracterização de uma variável aleatória discreta; distribuições de probabilidade: Uniforme, Bernoulli, Binomial, Poisson, Geométrica, Binomial Negativa e Hipergeométrica.5)Variáveis Aleatórias contínuas: Caracterização de uma variável aleatória contínua; distribuições de probabilidade: Uniforme, Exponencial e Normal.6)Aproximações: Aproximação das distribuições Binomial e Poisson pela distribuição Normal.7)Teorema Central do Limite: Distribuição da média amostral; distribuição da proporção amostral; intervalos de confiança para a média amostral e para a proporção amostral; dimensionamento amostral.8)Conceitos de Testes de Hipóteses: Erro Tipo I e Erro Tipo II; p-valor; poder.9)Testes de Hipóteses para uma única amostra: Teste de hipótese para a média; teste de hipótese para a proporção e teste de hipótese para a variância.10)Testes de Hipóteses para comparação de duas amostras: Teste de hipótese para comparação de médias (amostras independentes e dependentes); teste de hipótese para comparação de duas proporções e teste de hipótese para comparação de variâncias.11) Análise de Variância: Estimação do modelo; tabela de análise de variância; intervalos de confiança para a diferença entre as médias; correção de Bonferroni; teste de homocedasticidade.12)Regressão Linear Simples e Regressão Linear Múltipla: Estimação do modelo; interpretação dos parâmetros; tabela de análise de variância; intervalos de confiança para os parâmetros; R^2; análise dos resíduos."
$syllabusEn = "1)Descriptive Statistics: Population and sample; graphical presentation of data in statistics; measures of central tendency position and dispersion.2)Sampling methods: Simple random sampling with replacement and simple random sampling without replacement.3)Introduction to probability: Probability concepts; events probability; conditional probability; independence; Bayes Theorem.4)Discrete Random Variables: Discrete variables characterization; probability distributions: Uniform, Bernoulli, Binomial, Poisson, Geometric, Negative Binomial and Hipergeometric.5)Continuous Random Variables: Continuous variables characterization; probability distributions: Uniform, Exponential and Normal.6)Approximations: Approximation of Binomial and Poisson distributions by Normal distribution.7)Central Limit Theorem: Distribution of sample mean; distribution of sample proportion; confidence intervals for means and proportion estimated; sample sizing.8)Hypothesis test concept:  Type I Error and Type II Error; p-value; power.9)Hypothesis test for a single sample: Hypothesis test for mean; hypothesis test for proportion and hypothesis test for variance.10)Hypothesis test for two samples comparison: hypothesis test for two means comparison (dependents and independents samples); hypothesis test for two proportions comparison and hypothesis test for two variances comparison.11)Analysis of variance: Model estimation; Analysis of Variance table; confidence intervals for means difference; Bonferroni correction; homoscedasticity test.12)Simple linear regression and Multiple linear regression: Model estimation; parameters interpretation; Analysis of Variance table; confidence intervals for the parameters; R^2; residuals analysis."
$metodoNf = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$nfGe5 = "NF≥ 5,0."
$nfRc = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$bibliografia = "BUSSAB, Wilton O., MORETTIN, Pedro A. Estatística básica. 5. Ed. São Paulo: Saraiva, 2006.`nDEVORE, Jay L Probabilidade e estatística para engenharia. São Paulo: Ed Thomson Pioneira, 2006.`nJOHNSON, Richard A.; WICHERN, Dean W. Applied multivariate statistical analysis. 5. ed. Upper- Saddle River: Prentice Hall, 2002.`nLARSON, Ron ; FARBER, Betsy.  Estatística aplicada. São Paulo. Ed. Prentice Hall Brasil, 2010.`nHOFFMANN, R. Estatística para economistas. 4. ed. São Paulo: Pioneira, 2006.`nRYAN, Thomas. Estatística moderna para engenharia. São Paulo: Ed. Campus, 2009.`nRUNGER, George C.; MONTGOMERY, Douglas. Estatística aplicada e probabilidade para engenheiros. São Paulo: Ed. LTC, 2009."

# --- Insert a new row at 13 for "Docentes responsaveis:" content, shifting rows 13-23 down to 14-24 ---
$ws.Rows("13:13").Insert()

# The inserted row copies formatting from the row above (A12); clear it so A13 carries no style,
# matching the target layout where row 13 only has B/C content (no A13 label).
$ws.Range("A13").Style = "Normal"

# --- Fix the column definitions: split the merged A:B style range into just column A ---
$ws.Range("B:B").ColumnWidth = 60.7109375

# --- Row 10 (Objetivos:) - add the Portuguese course-objectives paragraph ---
$ws.Range("B10").Value = $objetivosPt
$ws.Range("C10").Value = $objetivosPt

# --- Row 13 (new) - Docentes responsaveis content ---
$ws.Range("B13").Value = $docenteNome
$ws.Range("C13").Value = $docenteNome

# --- Row 14 (Programa resumido:) - short PT syllabus summary ---
$ws.Range("B14").Value = $programaResumidoPt
$ws.Range("C14").Value = $programaResumidoPt

# --- Row 15 (Short syllabus:) - short EN syllabus summary ---
$ws.Range("B15").Value = $shortSyllabusEn
$ws.Range("C15").Value = $shortSyllabusEn

# --- Row 16 (Programa:) - full PT syllabus ---
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# --- Row 17 (Syllabus:) - full EN syllabus; add row height to match target (120pt) ---
$ws.Range("B17").Value = $syllabusEn
$ws.Range("C17").Value = $syllabusEn

# --- Row 19 (Metodo:) - evaluation method text ---
$ws.Range("B19").Value = $metodoNf
$ws.Range("C19").Value = $metodoNf

# --- Row 20 (Criterio:) - passing grade criterion ---
$ws.Range("B20").Value = $nfGe5
$ws.Range("C20").Value = $nfGe5

# --- Row 21 (Norma de recuperacao:) - recovery norm ---
$ws.Range("B21").Value = $nfRc
$ws.Range("C21").Value = $nfRc

# --- Row 22 (Bibliografia:) - bibliography list ---
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

